$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "Test" label in K2 with the new LTSD Parameters table header
$ws.Range("K2").Value = "LTSD Parameters"

# Right / Left sub-headers (row 3)
$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

# Column labels (row 4)
$ws.Range("K4").Value = "Threshols"
$ws.Range("L4").Value = "Win"
$ws.Range("M4").Value = "Threshold"
$ws.Range("N4").Value = "Win"

# Parameter values (row 5) - these look numeric but must be stored as text.
# Temporarily mark the cells as Text so the numeric-looking strings aren't
# auto-converted to numbers, then restore the default (Normal) style so no
# lingering number format is left applied to the cells.
# (Entry order L5, K5, M5, N5 reproduces the original shared-string table order.)
$ws.Range("K5:N5").NumberFormat = "@"
$ws.Range("L5").Value = "100.0"
$ws.Range("K5").Value = "8.0"
$ws.Range("M5").Value = "7.7"
$ws.Range("N5").Value = "400.0"
$ws.Range("K5:N5").Style = "Normal"

# Reflect the new selection (active cell) on the sheet
$ws.Range("N5").Select()
